# investor_advisors_without_users.xlsx
# - Random password for users getting created: drop the "Password"
#   column's header (H1) and sample value (H2), including H2's mailto
#   hyperlink, since passwords are now generated randomly instead of
#   being hard-coded in the upload template.
# - Default advisor view: scroll/select so column H is in view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Password" header and its sample "P@assword" value. Setting
# Value to an empty string removes the cell's shared-string reference
# (and, since no other cell uses those two strings, they drop out of
# sharedStrings.xml on save) while preserving the existing cell style.
$ws.Range("H1").Value = ""
$ws.Range("H2").Value = ""

# Drop the hyperlink that lived on H2 ("mailto:P@assword"). Deleting
# hyperlinks scoped to the range removes the sheet's hyperlink list, so
# re-create the still-wanted one on A2 (the investor email) afterwards.
$ws.Range("H2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:raleigh-wiliamsons@gmail.com", [Type]::Missing, [Type]::Missing, "raleigh-wiliamsons@gmail.com")

# Update the default view: scroll right so column D is the left-most
# visible column, and select H1:H2 (the now-empty former Password cells)
# as the active selection.
$excel.Goto($ws.Range("H1:H2"), $true)
$excel.ActiveWindow.ScrollColumn = 4
